$wb = $excel.ActiveWorkbook

# "Spain" is the last existing market sheet and serves as the template for
# the new "Netherlands" market sheet - duplicate it, placing the copy
# immediately after the original.
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy([System.Reflection.Missing]::Value, $spain)

# Excel makes the freshly-copied sheet the active sheet.
$netherlands = $wb.ActiveSheet
$netherlands.Name = "Netherlands"

# Fill in the market-specific values for the new sheet.
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Range("B4").Value = "NGC-3144/T2176 "

# Row 2 (wrapped "... Market" text) grows a touch to fit two lines, while
# rows 3-5 (previously tall to fit the Spain part number) go back to the
# sheet's standard height.
$netherlands.Rows(2).RowHeight = 28.8
$netherlands.Rows(3).AutoFit()
$netherlands.Rows(4).AutoFit()
$netherlands.Rows(5).AutoFit()

# The now-inactive "Spain" sheet reverts to a whole-sheet selection, and the
# new "Netherlands" sheet becomes the active tab selected at B4 (the order
# of these two Select() calls matters: the last one wins the active tab).
$spain.Range("A1:XFD1048576").Select() | Out-Null
$netherlands.Range("B4").Select() | Out-Null
